$wb = $excel.ActiveWorkbook

# ALC row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 789.2083
$ws.Range("J129").Value = 968.64703
$ws.Range("L129").Value = 2905.94109
$ws.Range("N129").Value = -12905.94109

# ALC row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 3469.52
$ws.Range("I135").Value = 915
$ws.Range("J135").Value = 10038.286
$ws.Range("K135").Value = 8235
$ws.Range("L135").Value = 90344.57399999999
$ws.Range("M135").Value = -5700
$ws.Range("N135").Value = -95414.57399999999

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 6017.273
$ws.Range("I141").Value = 6017.273
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 18051.819
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -12871.819
$ws.Range("N141").ClearContents()

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 693.42
$ws.Range("I2").Value = 597.125
$ws.Range("J2").Value = 1078.6
$ws.Range("K2").Value = 597.125
$ws.Range("L2").Value = 1078.6
$ws.Range("M2").Value = -484.125
$ws.Range("N2").Value = -1304.6

# ARM row 18
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H18").Value = 60000
$ws.Range("J18").Value = 60000
$ws.Range("L18").Value = 60000
$ws.Range("N18").Value = -60644

# ARM row 26
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 1244.1666
$ws.Range("I26").Value = 1244.1666
$ws.Range("K26").Value = 1244.1666
$ws.Range("M26").Value = -914.1666

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14233.27
$ws.Range("I32").Value = 4118.701
$ws.Range("J32").Value = 81923.08
$ws.Range("K32").Value = 4118.701
$ws.Range("L32").Value = 81923.08
$ws.Range("M32").Value = -3831.701
$ws.Range("N32").Value = -82497.08

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1146.6666
$ws.Range("I74").Value = 822.0625
$ws.Range("J74").Value = 2630.5715
$ws.Range("K74").Value = 822.0625
$ws.Range("L74").Value = 2630.5715
$ws.Range("M74").Value = 51.9375
$ws.Range("N74").Value = -4378.5715

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1146.6666
$ws.Range("I77").Value = 822.0625
$ws.Range("J77").Value = 2630.5715
$ws.Range("K77").Value = 4110.3125
$ws.Range("L77").Value = 13152.8575
$ws.Range("M77").Value = 257.6875
$ws.Range("N77").Value = -21888.8575

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 693.42
$ws.Range("I116").Value = 597.125
$ws.Range("J116").Value = 1078.6
$ws.Range("K116").Value = 597.125
$ws.Range("L116").Value = 1078.6
$ws.Range("M116").Value = 1696.875
$ws.Range("N116").Value = -5666.6

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 693.42
$ws.Range("I3").Value = 597.125
$ws.Range("J3").Value = 1078.6
$ws.Range("K3").Value = 597.125
$ws.Range("L3").Value = 1078.6
$ws.Range("M3").Value = -483.125
$ws.Range("N3").Value = -1306.6

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1972.1
$ws.Range("I58").Value = 2206.5217
$ws.Range("J58").Value = 1654.9412
$ws.Range("K58").Value = 2206.5217
$ws.Range("L58").Value = 1654.9412
$ws.Range("M58").Value = -2003.5217
$ws.Range("N58").Value = -2060.9412

# CRP row 94
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 369808.4
$ws.Range("I94").Value = 334210.34
$ws.Range("J94").Value = 386238.3
$ws.Range("K94").Value = 334210.34
$ws.Range("L94").Value = 386238.3
$ws.Range("M94").Value = -333759.34
$ws.Range("N94").Value = -387140.3

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1708.6207
$ws.Range("I132").Value = 1329.5264
$ws.Range("J132").Value = 2428.9
$ws.Range("K132").Value = 3988.5792
$ws.Range("L132").Value = 7286.700000000001
$ws.Range("M132").Value = -1458.5792
$ws.Range("N132").Value = -12346.7

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1929.65
$ws.Range("I134").Value = 2243.8
$ws.Range("J134").Value = 987.2
$ws.Range("K134").Value = 6731.400000000001
$ws.Range("L134").Value = 2961.6
$ws.Range("M134").Value = -4196.400000000001
$ws.Range("N134").Value = -8031.6

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1972.1
$ws.Range("I136").Value = 2206.5217
$ws.Range("J136").Value = 1654.9412
$ws.Range("K136").Value = 6619.5651
$ws.Range("L136").Value = 4964.8236
$ws.Range("M136").Value = -4069.5651
$ws.Range("N136").Value = -10064.8236

# CUL row 34
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1017.61536
$ws.Range("I34").Value = 1264.5
$ws.Range("J34").Value = 907.8889
$ws.Range("K34").Value = 3793.5
$ws.Range("L34").Value = 2723.6667
$ws.Range("M34").Value = -3709.5
$ws.Range("N34").Value = -2891.6667

# CUL row 39
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 1361.6
$ws.Range("J39").Value = 1457.3334
$ws.Range("L39").Value = 4372.0002
$ws.Range("N39").Value = -4960.0002

# CUL row 55
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 3140
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 3140
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 9420
$ws.Range("N55").Value = -9774
$ws.Range("M55").ClearContents()

# CUL row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 16743.125
$ws.Range("I68").Value = 718.4
$ws.Range("J68").Value = 43451
$ws.Range("K68").Value = 2155.2
$ws.Range("L68").Value = 130353
$ws.Range("M68").Value = -1344.2
$ws.Range("N68").Value = -131975

# CUL row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 16743.125
$ws.Range("I71").Value = 718.4
$ws.Range("J71").Value = 43451
$ws.Range("K71").Value = 6465.599999999999
$ws.Range("L71").Value = 391059
$ws.Range("M71").Value = -2409.599999999999
$ws.Range("N71").Value = -399171

# CUL row 81
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 2619.2856
$ws.Range("I81").Value = 1000
$ws.Range("J81").Value = 3833.75
$ws.Range("K81").Value = 3000
$ws.Range("L81").Value = 11501.25
$ws.Range("M81").Value = -1877
$ws.Range("N81").Value = -13747.25

# CUL row 84
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H84").Value = 2619.2856
$ws.Range("I84").Value = 1000
$ws.Range("J84").Value = 3833.75
$ws.Range("K84").Value = 9000
$ws.Range("L84").Value = 34503.75
$ws.Range("M84").Value = -3384
$ws.Range("N84").Value = -45735.75

# CUL row 97
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 632.6818
$ws.Range("I97").Value = 215
$ws.Range("J97").Value = 1050.3636
$ws.Range("K97").Value = 645
$ws.Range("L97").Value = 3151.0908
$ws.Range("M97").Value = -149
$ws.Range("N97").Value = -4143.0908

# CUL row 98
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 1050.3478
$ws.Range("I98").Value = 950
$ws.Range("J98").Value = 1142.3334
$ws.Range("K98").Value = 2850
$ws.Range("L98").Value = 3427.0002
$ws.Range("M98").Value = -1352
$ws.Range("N98").Value = -6423.0002

# CUL row 107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1295.1154
$ws.Range("I107").Value = 439.2857
$ws.Range("J107").Value = 1610.421
$ws.Range("K107").Value = 1317.8571
$ws.Range("L107").Value = 4831.263
$ws.Range("M107").Value = 602.1428999999998
$ws.Range("N107").Value = -8671.262999999999

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 907.30206
$ws.Range("I131").Value = 330
$ws.Range("J131").Value = 913.37897
$ws.Range("K131").Value = 990
$ws.Range("L131").Value = 2740.13691
$ws.Range("M131").Value = 4050
$ws.Range("N131").Value = -12820.13691

# GSM row 17
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 22000
$ws.Range("I17").Value = 500
$ws.Range("J17").Value = 27375
$ws.Range("K17").Value = 500
$ws.Range("L17").Value = 27375
$ws.Range("M17").Value = -332
$ws.Range("N17").Value = -27711

# GSM row 28
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H28").Value = 14932
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 14932
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 14932
$ws.Range("N28").Value = -15316
$ws.Range("M28").ClearContents()

# GSM row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2000
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 2000
$ws.Range("N113").Value = -6340
$ws.Range("M113").ClearContents()

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1704.4
$ws.Range("I126").Value = 1683.2858
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 5049.857400000001
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -2579.857400000001
$ws.Range("N126").Value = -10940

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2889.85
$ws.Range("I132").Value = 2551.4243
$ws.Range("J132").Value = 4485.2856
$ws.Range("K132").Value = 7654.2729
$ws.Range("L132").Value = 13455.8568
$ws.Range("M132").Value = -5124.2729
$ws.Range("N132").Value = -18515.8568

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 10000
$ws.Range("I61").Value = 10000
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 10000
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -9798
$ws.Range("N61").ClearContents()

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 10000
$ws.Range("I113").Value = 10000
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 10000
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -7830
$ws.Range("N113").ClearContents()

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 6921.972
$ws.Range("I132").Value = 4290.7827
$ws.Range("J132").Value = 11577.154
$ws.Range("K132").Value = 12872.3481
$ws.Range("L132").Value = 34731.462
$ws.Range("M132").Value = -10342.3481
$ws.Range("N132").Value = -39791.462

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 706.2222
$ws.Range("I136").Value = 581.7917
$ws.Range("J136").Value = 1701.6666
$ws.Range("K136").Value = 1745.3751
$ws.Range("L136").Value = 5104.9998
$ws.Range("M136").Value = 804.6249
$ws.Range("N136").Value = -10204.9998
